# Update column G ("K" - strikeouts) values on Sheet1 to reflect the
# regenerated save_data (K instead of Strike#).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$updates = @{
    2  = 2
    3  = 1
    4  = 2
    5  = 5
    6  = 0
    7  = 2
    8  = 0
    9  = 2
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 2
    15 = 1
    16 = 2
    17 = 1
    18 = 1
    19 = 3
    20 = 1
    22 = 2
    23 = 2
    25 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
